$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 15: plain (no highlight) pair
$ws.Range("A15").Value = "Âm Long Trực"
$ws.Range("B15").Value = "Bạn là người thông minh, biết cách ứng xử phù hợp và nên giữ đức tính nhu thuận làm kim chỉ nam cuộc đời để gặp nhiều may mắn. Thuận thiên vô chiến tự nhiên thành."

# Add row 16: A16 uses the highlighted style (same as A1:A7, A13, A14), B16 plain
$ws.Range("A16").Value = "Tuế Hổ Phù"
$ws.Range("B16").Value = "Bạn sinh ra gặp rất nhiều sóng gió cuộc đời nhưng đến khi vào đại vận trung niên bạn sẽ được hưởng trọn vẹn thành quả của những cố gắng, kiến thức, trải nghiệm đã đựợc tích lũy trước đó, cuộc sống gắn liền phần nhiều đến tín ngưỡng và tôn giáo. "

# Copy the highlighted fill style from A14 (which already has it) onto A16
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the active selection to match target (B16 selected)
$ws.Range("B16").Select() | Out-Null

$wb.Save()
